$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "Meta description: ..." paragraph that currently sits
#    right under the title heading.
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# ------------------------------------------------------------------
# 2. Insert a new bold paragraph ("Play Epic MONOPOLY II for Free -
#    Exciting Gameplay and Features") right before the last paragraph
#    (the old "Prompt for DALLE: ..." paragraph).
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastPara.Range.InsertParagraphBefore()

$count2 = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($count2 - 1)

$newParaXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Epic MONOPOLY II for Free - Exciting Gameplay and Features</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara.Range.InsertXML($newParaXml)

# ------------------------------------------------------------------
# 3. Replace the old DALLE image-prompt text (now the very last
#    paragraph) with the meta-description sentence, keeping the
#    paragraph's existing italic run formatting intact.
# ------------------------------------------------------------------
$oldText = "Prompt for DALLE: Create a feature image that captures the essence of Epic MONOPOLY II. The image should be in cartoon style, featuring a happy Maya warrior wearing glasses. The warrior should be standing in front of a Monopoly board with an overjoyed expression on their face. The Monopoly board should be full of colorful houses, hotels, and tokens. There should also be a large " + [char]34 + "Epic Wheel" + [char]34 + " in the background, hinting at the game's exciting bonus feature. The image should have a vibrant and energetic feel, with bold colors that pop and catch the player's eye. The Maya warrior should be an eye-catching element in the foreground, drawing attention to the game's theme and gameplay. The overall design should be playful and fun, creating a sense of excitement and anticipation in the player."
$newText = "Experience the classic board game with Epic MONOPOLY II online slot. Play for free with exciting bonuses and a maximum payout of 1,500x your bet."

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
